{"js": "// Replace the 25 two-digit multiplication problems in the table with new\n// problems/answers, matching each old value positionally (row, col) inside\n// the single 20-row x 5-col table. Blank spacer rows are left untouched.\n\nconst replacements = [\n  [\"86\u00d723=1978\", \"59\u00d762=3658\"],\n  [\"20\u00d729=580\", \"71\u00d759=4189\"],\n  [\"96\u00d759=5664\", \"40\u00d718=720\"],\n  [\"89\u00d722=1958\", \"80\u00d736=2880\"],\n  [\"85\u00d750=4250\", \"58\u00d721=1218\"],\n  [\"87\u00d790=7830\", \"65\u00d764=4160\"],\n  [\"17\u00d730=510\", \"96\u00d764=6144\"],\n  [\"99\u00d758=5742\", \"19\u00d721=399\"],\n  [\"15\u00d780=1200\", \"49\u00d779=3871\"],\n  [\"98\u00d747=4606\", \"70\u00d778=5460\"],\n  [\"93\u00d799=9207\", \"60\u00d780=4800\"],\n  [\"97\u00d732=3104\", \"78\u00d714=1092\"],\n  [\"15\u00d795=1425\", \"49\u00d727=1323\"],\n  [\"26\u00d728=728\", \"40\u00d757=2280\"],\n  [\"53\u00d758=3074\", \"59\u00d789=5251\"],\n  [\"14\u00d737=518\", \"74\u00d759=4366\"],\n  [\"59\u00d786=5074\", \"81\u00d726=2106\"],\n  [\"26\u00d770=1820\", \"85\u00d784=7140\"],\n  [\"96\u00d722=2112\", \"54\u00d765=3510\"],\n  [\"64\u00d798=6272\", \"84\u00d776=6384\"],\n  [\"70\u00d747=3290\", \"41\u00d759=2419\"],\n  [\"99\u00d741=4059\", \"80\u00d793=7440\"],\n  [\"66\u00d790=5940\", \"96\u00d744=4224\"],\n  [\"32\u00d759=1888\", \"26\u00d754=1404\"],\n  [\"70\u00d779=5530\", \"67\u00d770=4690\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    cell.load(\"value\");\n  }\n  await context.sync();\n\n  for (const cell of cells.items) {\n    const current = cell.value;\n    const hit = replacements.find(([oldText]) => oldText === current);\n    if (hit) {\n      cell.value = hit[1];\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit multiplication problems in the table with new\n# problems/answers. Each old string is unique in the document, so a plain\n# Find/Replace-All per pair reproduces the diff exactly.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"86\u00d723=1978\", \"59\u00d762=3658\"),\n    @(\"20\u00d729=580\",  \"71\u00d759=4189\"),\n    @(\"96\u00d759=5664\", \"40\u00d718=720\"),\n    @(\"89\u00d722=1958\", \"80\u00d736=2880\"),\n    @(\"85\u00d750=4250\", \"58\u00d721=1218\"),\n    @(\"87\u00d790=7830\", \"65\u00d764=4160\"),\n    @(\"17\u00d730=510\",  \"96\u00d764=6144\"),\n    @(\"99\u00d758=5742\", \"19\u00d721=399\"),\n    @(\"15\u00d780=1200\", \"49\u00d779=3871\"),\n    @(\"98\u00d747=4606\", \"70\u00d778=5460\"),\n    @(\"93\u00d799=9207\", \"60\u00d780=4800\"),\n    @(\"97\u00d732=3104\", \"78\u00d714=1092\"),\n    @(\"15\u00d795=1425\", \"49\u00d727=1323\"),\n    @(\"26\u00d728=728\",  \"40\u00d757=2280\"),\n    @(\"53\u00d758=3074\", \"59\u00d789=5251\"),\n    @(\"14\u00d737=518\",  \"74\u00d759=4366\"),\n    @(\"59\u00d786=5074\", \"81\u00d726=2106\"),\n    @(\"26\u00d770=1820\", \"85\u00d784=7140\"),\n    @(\"96\u00d722=2112\", \"54\u00d765=3510\"),\n    @(\"64\u00d798=6272\", \"84\u00d776=6384\"),\n    @(\"70\u00d747=3290\", \"41\u00d759=2419\"),\n    @(\"99\u00d741=4059\", \"80\u00d793=7440\"),\n    @(\"66\u00d790=5940\", \"96\u00d744=4224\"),\n    @(\"32\u00d759=1888\", \"26\u00d754=1404\"),\n    @(\"70\u00d779=5530\", \"67\u00d770=4690\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]2) | Out-Null\n}\n"}
